# ------------------------------------------------------------------
# PlayerPerformance_3586.xlsx edit:
#  1. Insert a new "Player Info" sheet at the front.
#  2. Keep "ODI Batting" / "ODI Bowling" in place (now shifted right by
#     one tab) and:
#       - rename column D header MATCH_CARD_LINK -> MATCH_CODE
#       - replace the full howstat URL with just the numeric match code
#       - drop the stray empty INNING_NUMBER cells (ODI Batting only)
#  3. Append a new "ODI Batting Extra" sheet at the end.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBowling = $wb.Worksheets.Item("ODI Bowling")

# ---- helper: strip "...MatchCode=1234" down to "1234" on a MATCH_CARD_LINK
#      column, renaming its header to MATCH_CODE -------------------------
function Convert-MatchLinkColumn($ws, $colIndex) {
    $headerCell = $ws.Cells.Item(1, $colIndex)
    $headerCell.NumberFormat = "@"
    $headerCell.Value = "MATCH_CODE"

    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $colIndex)
        $val = $cell.Value()
        if ($val -ne $null -and $val -match "MatchCode=(\d+)") {
            $cell.NumberFormat = "@"
            $cell.Value = $matches[1]
        }
    }
}

# ===== ODI Batting (column D = MATCH_CARD_LINK) =====
Convert-MatchLinkColumn $odiBatting 4

# Rows whose INNING_NUMBER (col B) is blank lost their stray empty cell.
foreach ($r in 52, 82, 103, 143, 145) {
    $cell = $odiBatting.Cells.Item($r, 2)
    if ($cell.Value() -eq $null -or $cell.Value() -eq "") {
        $cell.Value = ""
    }
}

# ===== ODI Bowling (column B = MATCH_CARD_LINK) =====
Convert-MatchLinkColumn $odiBowling 2

# ------------------------------------------------------------------
# New "ODI Batting Extra" sheet, appended after "ODI Bowling".
# (Must be added *before* the "Player Info" insert below -- inserting
# a sheet ahead of "ODI Batting"/"ODI Bowling" first throws off the
# "after odiBowling" positional add.)
# ------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Add($null, $odiBowling)
$battingExtra.Name = "ODI Batting Extra"

$beHeaders = New-Object 'object[,]' 1,6
$beHeaders[0,0] = "MATCH_CODE"
$beHeaders[0,1] = "BATTING_POSITION"
$beHeaders[0,2] = "NUM_4"
$beHeaders[0,3] = "NUM_6"
$beHeaders[0,4] = "PERCENT_RUNS_OF_TOTAL"
$beHeaders[0,5] = "MAN_OF_MATCH"
$beHeaderRange = $battingExtra.Range("A1:F1")
$beHeaderRange.NumberFormat = "@"
$beHeaderRange.Value = $beHeaders
$beHeaderRange.Font.Bold = $true
$beHeaderRange.HorizontalAlignment = -4108
$beHeaderRange.VerticalAlignment = -4160
$beHeaderRange.Borders.LineStyle = 1

# MATCH_CODE, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH are text;
# BATTING_POSITION is numeric.
$beRows = @(
    @("4466", "2", "8", "0", "35.57%", "NO"),
    @("4467", "2", "6", "0", "32.91%", "NO"),
    @("4468", "2", "1", "0", "3.07%",  "NO"),
    @("4474", "",  "",  "",  "",       "NO"),
    @("4475", "1", "4", "0", "9.31%",  "NO"),
    @("4478", "",  "",  "",  "",       "NO"),
    @("4492", "2", "3", "1", "14.04%", "NO"),
    @("4494", "",  "",  "",  "",       "NO"),
    @("4496", "2", "5", "2", "36.44%", "NO"),
    @("4520", "2", "3", "1", "12.50%", "NO"),
    @("4522", "",  "",  "",  "",       "NO"),
    @("4605", "1", "1", "0", "1.67%",  "NO"),
    @("4608", "1", "0", "0", "",       "NO"),
    @("4614", "1", "14", "5", "33.43%", "NO"),
    @("4693", "",  "",  "",  "",       "NO"),
    @("4694", "1", "4", "2", "15.31%", "NO"),
    @("4696", "",  "",  "",  "",       "NO"),
    @("4726", "",  "",  "",  "",       "NO"),
    @("4729", "",  "",  "",  "",       "NO"),
    @("4734", "2", "1", "0", "6.93%",  "NO")
)

$rowCount = $beRows.Count
$beTextCols = New-Object 'object[,]' $rowCount,6
for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $beRows[$i]
    $beTextCols[$i,0] = $row[0]
    $beTextCols[$i,2] = $row[2]
    $beTextCols[$i,3] = $row[3]
    $beTextCols[$i,4] = $row[4]
    $beTextCols[$i,5] = $row[5]
}

$lastRow = 1 + $rowCount
$beBodyRange = $battingExtra.Range($battingExtra.Cells.Item(2,1), $battingExtra.Cells.Item($lastRow,6))
$beBodyRange.NumberFormat = "@"
$beBodyRange.Value = $beTextCols

# BATTING_POSITION (column B) is numeric where present, blank otherwise.
for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $beRows[$i]
    $r = 2 + $i
    $posCell = $battingExtra.Cells.Item($r, 2)
    if ($row[1] -ne "") {
        $posCell.NumberFormat = "General"
        $posCell.Value = [double]$row[1]
    }
}

# ------------------------------------------------------------------
# New "Player Info" sheet, inserted before "ODI Batting".
# ------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

$piHeaders = New-Object 'object[,]' 1,4
$piHeaders[0,0] = "ID"
$piHeaders[0,1] = "NAME"
$piHeaders[0,2] = "BATTING_HAND"
$piHeaders[0,3] = "BOWL_STYLE"
$piHeaderRange = $playerInfo.Range("A1:D1")
$piHeaderRange.NumberFormat = "@"
$piHeaderRange.Value = $piHeaders
$piHeaderRange.Font.Bold = $true
$piHeaderRange.HorizontalAlignment = -4108
$piHeaderRange.VerticalAlignment = -4160
$piHeaderRange.Borders.LineStyle = 1

$piData = New-Object 'object[,]' 1,4
$piData[0,0] = "3586"
$piData[0,1] = "Paul Robert Stirling"
$piData[0,2] = "Right Handed"
$piData[0,3] = "Right Arm Off Break"
$piDataRange = $playerInfo.Range("A2:D2")
$piDataRange.NumberFormat = "@"
$piDataRange.Value = $piData
